# Update "想去人数" (column F) values on both the "展览" and "全部类型"
# worksheets to reflect regenerated site stats.
$wb = $excel.ActiveWorkbook

# Row -> New value for column F
$updates = @{
    7  = 1319
    8  = 1562
    10 = 431
    15 = 114
    17 = 318
    19 = 1763
    26 = 4245
    32 = 614
    34 = 310
    36 = 155
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
